$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.818993210792542
$ws.Range("B1").Value = 4.573708057403564
$ws.Range("C1").Value = 4.064336776733398
$ws.Range("D1").Value = 0.9051988124847412
$ws.Range("E1").Value = 0.4756152033805847
